# Apply statement_29.xlsx edits: new account holder, new card number,
# shifted statement period (Jan/Feb 2024 instead of Apr/May 2025), and a
# regenerated set of transactions (one fewer row than before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: account holder name + card number ---
$ws.Range("C2").Value = "Hartmut"

# The card number is all-digits, so a plain assignment would be
# auto-coerced to a number (and silently lose precision past 15
# significant digits). Force it to text with a leading apostrophe, then
# restore the original cell formatting (the text coercion mints a new
# "@"-formatted style) by pasting formats from an unrelated cell that
# still carries the original style.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("D2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance date ---
$ws.Range("D5").Value = "KONTOSTAND AM 23.01.2024"

# --- Transaction rows 6-10 (BELEG date, VALUTA date, BESCHREIBUNG, BETRAG) ---
$ws.Range("B6").Value = "27.01."
$ws.Range("C6").Value = "28.01."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "24,75-"

$ws.Range("B7").Value = "31.01."
$ws.Range("C7").Value = "01.02."
$ws.Range("D7").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E7").Value = "47,51-"

$ws.Range("B8").Value = "04.02."
$ws.Range("C8").Value = "05.02."
$ws.Range("D8").Value = "AMAZON.DE MKTPLC EU UCNZKE"
$ws.Range("E8").Value = "54,46-"

$ws.Range("B9").Value = "07.02."
$ws.Range("C9").Value = "08.02."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-37892561"
$ws.Range("E9").Value = "55,47-"

$ws.Range("B10").Value = "10.02."
$ws.Range("C10").Value = "11.02."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 86351543"
$ws.Range("E10").Value = "83,80-"

# --- Row 11 previously held a 6th transaction (PAYPAL NWJHSP); the
#     regenerated statement only has 5, so the row becomes blank. The
#     amount cell keeps right-alignment but also picks up vertical-center
#     + wrap (matching the blank placeholder style used elsewhere). ---
$ws.Range("B11:D11").ClearContents()
$ws.Range("E11").Value = $null
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# --- Closing balance + next billing date ---
$ws.Range("D12").Value = "KONTOSTAND AM 13.02.2024"
$ws.Range("E12").Value = "265,99-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 19.02.2024"
